$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing "FillPageVehicleData" row (old row 5),
# pushing it down to row 6, and making room for the new "Goto insurant page" row.
$ws.Rows("5:5").Insert()

# New row 7: "FillPageInsurantData" smoke-test row (written first so its
# shared string lands before the other new strings, matching authoring order).
$ws.Range("A7").Value = "102_VehicleInsuranceAutomobile_001_SmokeTest_FillPageInsurantData"

# New row 5: "Goto insurant page" action, selecting a radiobutton-style control.
$ws.Range("A5").Value = "Goto insurant page"
$ws.Range("D5").Value = "<SELECT>"

# New row 8: "Button Next" action.
$ws.Range("A8").Value = "Button Next from Page VehicleData"
$ws.Range("C8").Value = "Button Next"

# Fill in the remaining (already-existing) shared-string cells for the new rows.
$ws.Range("B5").Value = "<SET>"
$ws.Range("H5").Value = "<NOP>"

$ws.Range("B7").Value = "<SET>"
$ws.Range("D7").Value = "102_VehicleInsuranceAutomobile_001_SmokeTest_FillPage"
$ws.Range("H7").Value = "<NOP>"

$ws.Range("B8").Value = "<SET>"
$ws.Range("H8").Value = "<NOP>"

# Column D now holds long strings like column C, so widen it to match.
$ws.Columns("D").ColumnWidth = $ws.Columns("C").ColumnWidth()

# Move the selection to D1 (matches the saved sheet view).
$ws.Range("D1").Select()

# The illustration picture keeps its on-sheet size but slides down with the
# newly inserted row.
$shp = $ws.Shapes.Item(1)
$shp.Top = 155.40007874015748
